# Scheduled market-data refresh: updates currentAveragePrice / LevePrice /
# LeveProfit columns (H-N) on several leve rows across the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets with newly fetched values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 146.66667
$ws.Range("I2").Value = 175
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 175
$ws.Range("L2").Value = 90
$ws.Range("M2").Value = -62
$ws.Range("N2").Value = -316

$ws.Range("H53").Value = 141.75
$ws.Range("J53").Value = 141.75
$ws.Range("L53").Value = 141.75
$ws.Range("N53").Value = -1415.75

$ws.Range("H70").Value = 2778.5
$ws.Range("I70").Value = 2921.375
$ws.Range("K70").Value = 8764.125
$ws.Range("M70").Value = -8494.125

$ws.Range("H73").Value = 2778.5
$ws.Range("I73").Value = 2921.375
$ws.Range("K73").Value = 8764.125
$ws.Range("M73").Value = -7828.125

$ws.Range("H98").Value = 1925.24
$ws.Range("I98").Value = 856.7
$ws.Range("K98").Value = 856.7
$ws.Range("M98").Value = 641.3

$ws.Range("H122").Value = 1925.24
$ws.Range("I122").Value = 856.7
$ws.Range("K122").Value = 2570.1
$ws.Range("M122").Value = -120.1000000000004

$ws.Range("H132").Value = 2496.4695
$ws.Range("I132").Value = 1381.6
$ws.Range("K132").Value = 4144.799999999999
$ws.Range("M132").Value = -1614.799999999999

$ws.Range("H138").Value = 2414.6667
$ws.Range("I138").Value = 1943.091
$ws.Range("J138").Value = 2650.4546
$ws.Range("K138").Value = 5829.272999999999
$ws.Range("L138").Value = 7951.3638
$ws.Range("M138").Value = -689.2729999999992
$ws.Range("N138").Value = -18231.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7875.9546
$ws.Range("I45").Value = 10549.692
$ws.Range("K45").Value = 10549.692
$ws.Range("M45").Value = -10172.692

$ws.Range("H74").Value = 1893.2858
$ws.Range("I74").Value = 1893.2858
$ws.Range("K74").Value = 1893.2858
$ws.Range("M74").Value = -1019.2858

$ws.Range("H77").Value = 1893.2858
$ws.Range("I77").Value = 1893.2858
$ws.Range("K77").Value = 9466.429
$ws.Range("M77").Value = -5098.429

$ws.Range("H117").Value = 82500
$ws.Range("J117").Value = 82500
$ws.Range("L117").Value = 82500
$ws.Range("N117").Value = -91678

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3103.6667
$ws.Range("I22").Value = 2931.2727
$ws.Range("K22").Value = 2931.2727
$ws.Range("M22").Value = -2758.2727

$ws.Range("H64").Value = 842.9
$ws.Range("I64").Value = 668
$ws.Range("J64").Value = 1017.8
$ws.Range("K64").Value = 668
$ws.Range("L64").Value = 1017.8
$ws.Range("M64").Value = -443
$ws.Range("N64").Value = -1467.8

$ws.Range("H67").Value = 842.9
$ws.Range("I67").Value = 668
$ws.Range("J67").Value = 1017.8
$ws.Range("K67").Value = 668
$ws.Range("L67").Value = 1017.8
$ws.Range("M67").Value = 112
$ws.Range("N67").Value = -2577.8

$ws.Range("H86").Value = 2857.55
$ws.Range("I86").Value = 3234.4375
$ws.Range("J86").Value = 1350
$ws.Range("K86").Value = 3234.4375
$ws.Range("L86").Value = 1350
$ws.Range("M86").Value = -2111.4375
$ws.Range("N86").Value = -3596

$ws.Range("H89").Value = 2857.55
$ws.Range("I89").Value = 3234.4375
$ws.Range("J89").Value = 1350
$ws.Range("K89").Value = 16172.1875
$ws.Range("L89").Value = 6750
$ws.Range("M89").Value = -10556.1875
$ws.Range("N89").Value = -17982

$ws.Range("H141").Value = 64758.832
$ws.Range("J141").Value = 64758.832
$ws.Range("L141").Value = 64758.832
$ws.Range("N141").Value = -75118.83199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2907.5
$ws.Range("J16").Value = 2400
$ws.Range("L16").Value = 2400
$ws.Range("N16").Value = -2974

$ws.Range("H31").Value = 12851.2705
$ws.Range("I31").Value = 3848.0454
$ws.Range("K31").Value = 3848.0454
$ws.Range("M31").Value = -3553.0454

$ws.Range("H34").Value = 12851.2705
$ws.Range("I34").Value = 3848.0454
$ws.Range("K34").Value = 3848.0454
$ws.Range("M34").Value = -3646.0454

$ws.Range("H86").Value = 3606.3572
$ws.Range("I86").Value = 3437.5
$ws.Range("J86").Value = 3831.5
$ws.Range("K86").Value = 3437.5
$ws.Range("L86").Value = 3831.5
$ws.Range("M86").Value = -2314.5
$ws.Range("N86").Value = -6077.5

$ws.Range("H89").Value = 3606.3572
$ws.Range("I89").Value = 3437.5
$ws.Range("J89").Value = 3831.5
$ws.Range("K89").Value = 17187.5
$ws.Range("L89").Value = 19157.5
$ws.Range("M89").Value = -11571.5
$ws.Range("N89").Value = -30389.5

$ws.Range("H113").Value = 2907.5
$ws.Range("J113").Value = 2400
$ws.Range("L113").Value = 2400
$ws.Range("N113").Value = -6740

$ws.Range("H135").Value = 152780
$ws.Range("J135").Value = 152780
$ws.Range("L135").Value = 152780
$ws.Range("N135").Value = -162920

$ws.Range("H141").Value = 300777
$ws.Range("J141").Value = 300777
$ws.Range("L141").Value = 300777
$ws.Range("N141").Value = -311137

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 394
$ws.Range("I2").Value = 1132.8
$ws.Range("J2").Value = 58.18182
$ws.Range("K2").Value = 6796.799999999999
$ws.Range("L2").Value = 349.09092
$ws.Range("M2").Value = -6683.799999999999
$ws.Range("N2").Value = -575.09092

$ws.Range("H34").Value = 892.7857
$ws.Range("I34").Value = 194.66667
$ws.Range("J34").Value = 1416.375
$ws.Range("K34").Value = 584.00001
$ws.Range("L34").Value = 4249.125
$ws.Range("M34").Value = -500.00001
$ws.Range("N34").Value = -4417.125

$ws.Range("H46").Value = 4399.75
$ws.Range("J46").Value = 9999
$ws.Range("L46").Value = 29997
$ws.Range("N46").Value = -30179

$ws.Range("H113").Value = 1734.1818
$ws.Range("J113").Value = 1873.8948
$ws.Range("L113").Value = 5621.6844
$ws.Range("N113").Value = -9961.6844

$ws.Range("H131").Value = 1701.5294
$ws.Range("J131").Value = 1723
$ws.Range("L131").Value = 5169
$ws.Range("N131").Value = -15249

$ws.Range("H139").Value = 4463.8125
$ws.Range("I139").Value = 4139.25
$ws.Range("K139").Value = 12417.75
$ws.Range("M139").Value = -7277.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3388.4443
$ws.Range("I80").Value = 2240
$ws.Range("J80").Value = 4824
$ws.Range("K80").Value = 2240
$ws.Range("L80").Value = 4824
$ws.Range("M80").Value = -1242
$ws.Range("N80").Value = -6820

$ws.Range("H83").Value = 3388.4443
$ws.Range("I83").Value = 2240
$ws.Range("J83").Value = 4824
$ws.Range("K83").Value = 11200
$ws.Range("L83").Value = 24120
$ws.Range("M83").Value = -6208
$ws.Range("N83").Value = -34104

$ws.Range("H107").Value = 6076.727
$ws.Range("I107").Value = 340.83334
$ws.Range("J107").Value = 12959.8
$ws.Range("K107").Value = 340.83334
$ws.Range("L107").Value = 12959.8
$ws.Range("M107").Value = 1579.16666
$ws.Range("N107").Value = -16799.8

$ws.Range("H113").Value = 3399.4583
$ws.Range("I113").Value = 2926.8
$ws.Range("K113").Value = 2926.8
$ws.Range("M113").Value = -756.8000000000002

$ws.Range("H122").Value = 2322.6667
$ws.Range("I122").Value = 2322.6667
$ws.Range("K122").Value = 6968.000100000001
$ws.Range("M122").Value = -4518.000100000001

$ws.Range("H126").Value = 1544.25
$ws.Range("I126").Value = 1100
$ws.Range("K126").Value = 3300
$ws.Range("M126").Value = -830

$ws.Range("H132").Value = 8978
$ws.Range("I132").Value = 8978
$ws.Range("K132").Value = 26934
$ws.Range("M132").Value = -24404

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 700
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 700
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H40").Value = 22249
$ws.Range("I40").Value = 37001.332
$ws.Range("J40").Value = 7496.6665
$ws.Range("K40").Value = 37001.332
$ws.Range("L40").Value = 7496.6665
$ws.Range("M40").Value = -36865.332
$ws.Range("N40").Value = -7768.6665

$ws.Range("H55").Value = 591.63635
$ws.Range("I55").Value = 751.8333
$ws.Range("J55").Value = 399.4
$ws.Range("K55").Value = 751.8333
$ws.Range("L55").Value = 399.4
$ws.Range("M55").Value = -578.8333
$ws.Range("N55").Value = -745.4

$ws.Range("H61").Value = 85863.875
$ws.Range("I61").Value = 78883.84
$ws.Range("K61").Value = 78883.84
$ws.Range("M61").Value = -78681.84

$ws.Range("H113").Value = 85863.875
$ws.Range("I113").Value = 78883.84
$ws.Range("K113").Value = 78883.84
$ws.Range("M113").Value = -76713.84

$ws.Range("H122").Value = 291285.56
$ws.Range("I122").Value = 670001.3
$ws.Range("J122").Value = 7248.75
$ws.Range("K122").Value = 2010003.9
$ws.Range("L122").Value = 21746.25
$ws.Range("M122").Value = -2007553.9
$ws.Range("N122").Value = -26646.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H122").Value = 1386.5
$ws.Range("J122").Value = 1137.6666
$ws.Range("L122").Value = 3412.9998
$ws.Range("N122").Value = -8312.9998

$ws.Range("H126").Value = 1700
$ws.Range("I126").Value = 1656.5625
$ws.Range("K126").Value = 4969.6875
$ws.Range("M126").Value = -2499.6875

$ws.Range("H132").Value = 1874.1613
$ws.Range("I132").Value = 1290.9524
$ws.Range("K132").Value = 3872.857199999999
$ws.Range("M132").Value = -1342.857199999999
